$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.105.22"
$ws.Range("E2").Value = "  -3.65%  "
$ws.Range("D3").Value = "3.291.60"
$ws.Range("E3").Value = "  -5.80%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'540.11"
$ws.Range("E5").Value = "  -2.67%  "
$ws.Range("D6").Value = "'170.03"
$ws.Range("E6").Value = "  -4.93%  "
$ws.Range("D7").Value = "'0.607"
$ws.Range("E7").Value = "  -4.74%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "3.283.13"
$ws.Range("E9").Value = "  -5.96%  "
$ws.Range("D10").Value = "'0.607"
$ws.Range("E10").Value = "  -3.65%  "
$ws.Range("D11").Value = "'0.151"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").Value = "'52.16"
$ws.Range("E12").Value = "  -2.85%  "
$ws.Range("D13").Value = "'0.0000261"
$ws.Range("E13").Value = "  -3.86%  "
$ws.Range("D14").Value = "'8.80"
$ws.Range("E14").Value = "  -4.73%  "
$ws.Range("D15").Value = "3.804.77"
$ws.Range("E15").Value = "  -6.21%  "
$ws.Range("D16").Value = "'17.95"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("E17").Value = "  -4.04%  "
$ws.Range("D18").Value = "3.273.67"
$ws.Range("E18").Value = "  -6.47%  "
$ws.Range("D19").Value = "'11.55"
$ws.Range("E19").Value = "  -4.73%  "
$ws.Range("D20").Value = "62.964.61"
$ws.Range("E20").Value = "  -3.96%  "
$ws.Range("D21").Value = "'0.961"
$ws.Range("E21").Value = "  -3.39%  "
$ws.Range("D22").Value = "'415.19"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").Value = "'4.36"
$ws.Range("E23").Value = "  +6.23%  "
$ws.Range("D24").Value = "'3.99"
$ws.Range("E24").Value = "  -1.29%  "
$ws.Range("D25").Value = "'13.33"
$ws.Range("E25").Value = "  +4.57%  "
$ws.Range("D26").Value = "'82.45"
$ws.Range("E26").Value = "  -4.05%  "
$ws.Range("D27").Value = "'10.53"
$ws.Range("E27").Value = "  -2.43%  "
$ws.Range("D28").Value = "'2.70"
$ws.Range("E28").Value = "  -5.06%  "
$ws.Range("D29").Value = "'8.51"
$ws.Range("E29").Value = "  -5.62%  "
$ws.Range("D30").Value = "'28.84"
$ws.Range("E30").Value = "  -4.67%  "
$ws.Range("D31").Value = "'6.32"
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("D32").Value = "'11.24"
$ws.Range("E32").Value = "  -3.50%  "
$ws.Range("D33").Value = "'568.59"
$ws.Range("E33").Value = "  -6.47%  "
$ws.Range("E34").Value = "  -3.99%  "
$ws.Range("D35").Value = "'57.59"
$ws.Range("E35").Value = "  -3.40%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").Value = "'0.145"
$ws.Range("E37").Value = "  -1.46%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "'3.40"
$ws.Range("E38").Value = "  +4.46%  "
$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").Value = "'34.86"
$ws.Range("E39").Value = "  -6.34%  "
$ws.Range("D40").Value = "0.0₃0733"
$ws.Range("E40").Value = "  -7.09%  "
$ws.Range("D41").Value = "'0.361"
$ws.Range("E41").Value = "  -4.86%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.090.45"
$ws.Range("E42").Value = "  -8.17%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'0.998"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("D44").Value = "'3.22"
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("D45").Value = "'2.72"
$ws.Range("E45").Value = "  -4.29%  "
$ws.Range("D46").Value = "'0.0397"
$ws.Range("E46").Value = "  -4.31%  "
$ws.Range("D47").Value = "'2.40"
$ws.Range("E47").Value = "  -5.16%  "
$ws.Range("D48").Value = "'2.57"
$ws.Range("E48").Value = "  -4.98%  "
$ws.Range("D49").Value = "'0.127"
$ws.Range("E49").Value = "  -3.75%  "
$ws.Range("D50").Value = "'131.99"
$ws.Range("E50").Value = "  -4.07%  "
$ws.Range("D51").Value = "'7.96"
$ws.Range("E51").Value = "  -5.64%  "
